# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-10 15:17:01
#
# Updates the "Recorded By" lists (reordered / refreshed), the related
# student-count and percentage figures on the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ANATOMY, session 1) - reorder "Recorded By" list
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# Row 3 (ANATOMY, session 2) - reorder + add recorder, bump attendance count
$ws.Range("G3").Value = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H3").Value = "89/251"

# Row 9 (HISTOLOGY, session 1) - reorder "Recorded By" list
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Row 10 - Average Attendance % stat, updated to match new figures.
# Force Text format first so Excel stores the literal "28.0%" string
# instead of auto-converting it to a percentage number.
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "28.0%"

# Row 15 - Average Attendance % (group statistics mirror of L10)
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "28.0%"

# Row 28 (PHYSIOLOGY, session 1) - reorder "Recorded By" list
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
